# Updates cryptos list prices/volumes per the scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.480.92"
$ws.Range("E2").Value = "  -2.89%  "

# Row 3
$ws.Range("D3").Value = "1.799.66"
$ws.Range("E3").Value = "  -2.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("D5").Value = "'229.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

# Row 6
$ws.Range("D6").Value = "'0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.41%  "

# Row 8
$ws.Range("D8").Value = "'39.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.85%  "

# Row 9
$ws.Range("E9").Value = "  +2.80%  "

# Row 10
$ws.Range("D10").Value = "'0.0678"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.93%  "

# Row 11
$ws.Range("D11").Value = "'0.0989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.12%  "

# Row 12
$ws.Range("D12").Value = "2.059.25"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13
$ws.Range("D13").Value = "'11.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.86%  "

# Row 14
$ws.Range("D14").Value = "1.794.90"
$ws.Range("E14").Value = "  -2.50%  "

# Row 15
$ws.Range("D15").Value = "'0.657"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "

# Row 16
$ws.Range("D16").Value = "'4.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.81%  "

# Row 17
$ws.Range("D17").Value = "34.329.05"
$ws.Range("E17").Value = "  -3.29%  "

# Row 18
$ws.Range("D18").Value = "'68.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.15%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0778"
$ws.Range("E19").Value = "  -2.95%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'239.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.01%  "

# Row 21
$ws.Range("D21").Value = "'11.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "

# Row 22
$ws.Range("D22").Value = "'4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "

# Row 23
$ws.Range("E23").Value = "  +0.43%  "

# Row 24
$ws.Range("E24").Value = "  -0.76%  "

# Row 25
$ws.Range("D25").Value = "'173.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

# Row 26
$ws.Range("D26").Value = "'7.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "

# Row 27
$ws.Range("D27").Value = "'17.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.66%  "

# Row 28
$ws.Range("D28").Value = "'0.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29
$ws.Range("D29").Value = "'1.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.23%  "

# Row 30
$ws.Range("E30").Value = "  +0.34%  "

# Row 31
$ws.Range("E31").Value = "  +1.19%  "

# Row 32
$ws.Range("D32").Value = "'0.0541"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "

# Row 33
$ws.Range("D33").Value = "'3.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.01%  "

# Row 34
$ws.Range("E34").Value = "  +7.66%  "

# Row 35
$ws.Range("E35").Value = "  -2.79%  "

# Row 36
$ws.Range("D36").Value = "'0.694"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37
$ws.Range("D37").Value = "'90.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.05%  "

# Row 38
$ws.Range("E38").Value = "  +4.91%  "

# Row 39
$ws.Range("D39").Value = "1.323.00"
$ws.Range("E39").Value = "  -1.79%  "

# Row 40
$ws.Range("E40").Value = "  -3.08%  "

# Row 41
$ws.Range("D41").Value = "'0.957"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.72%  "

# Row 42
$ws.Range("D42").Value = "'2.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.85%  "

# Row 43
$ws.Range("D43").Value = "'14.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.19%  "

# Row 44
$ws.Range("E44").Value = "  -9.54%  "

# Row 45
$ws.Range("E45").Value = "  -3.68%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'6.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0512"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "

# Row 48
$ws.Range("D48").Value = "1.981.44"
$ws.Range("E48").Value = "  -1.67%  "

# Row 49
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$ws.Range("E50").Value = "  +3.83%  "

# Row 51
$ws.Range("D51").Value = "'97.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.83%  "

